# The 'reviews_count' column (column E) is removed from the worksheet.
# All columns to the right (reviews_average, latitude, longitude,
# is_permanently_closed, gmaps_link, latest_review_date) shift one
# position to the left (F->E, G->F, H->G, I->H, J->I, K->J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(5).Delete()
